$d = $word.ActiveDocument

# Locate the résumé line that currently reads "性别：男" (gender: male).
$genderPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*性别：男*") {
        $genderPara = $p
        break
    }
}

if ($genderPara -ne $null) {
    # The edit amounts to: put the cursor at the end of the previous line
    # (the "姓名：..." line), press Enter - which mints a new paragraph mark
    # that inherits the preceding run's formatting (rFonts hint="eastAsia") -
    # retype "性别：男" into that new line, and then type the new line of text
    # ("毕业于湖南信息职业技术学院") into what used to be the "性别：男"
    # paragraph (whose own paragraph mark / bookmark stay untouched).
    $prevPara = $genderPara.Previous()
    $prevPara.Range.InsertParagraphAfter()

    $newPara = $prevPara.Next()
    $newPara.Range.Text = "性别：男"

    # The original "性别：男" paragraph (with its bookmark) is now the one
    # after the freshly typed line; only its run text changes.
    $target = $newPara.Next()
    $target.Range.Find.Execute("性别：男", $true, $false, $false, $false, $false, `
                                $true, 1, $false, "毕业于湖南信息职业技术学院", 1)
}
